# Add a new "Sampling" column (derived from Season) to table_sample_stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; everything from J onward (Tissue_water, Location,
# final/raw, Comment, ! Exchange) shifts one column to the right (K..O).
$ws.Columns("J:J").Insert()

# Header for the new column.
$ws.Range("J2").Value = "Sampling"

# Populate the new "Sampling" column based on the existing "Season" column (I):
#   Winter -> S1, Spring -> S2, no -> no
for ($r = 3; $r -le 20; $r++) {
    $ws.Cells.Item($r, 10).Value = "S1"
}
for ($r = 21; $r -le 32; $r++) {
    $ws.Cells.Item($r, 10).Value = "S2"
}
for ($r = 33; $r -le 40; $r++) {
    $ws.Cells.Item($r, 10).Value = "no"
}

# Re-apply the AutoFilter so its stored range grows from A2:M40 to A2:N40.
[void]$ws.Range("A2:N40").AutoFilter()
[void]$ws.Range("A2:N40").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new filter range.
$name = $wb.Names.Item("_xlnm._FilterDatabase")
$name.RefersTo = "=table_sample_stats!`$A`$2:`$N`$40"

# Update the selected cell to reflect where editing ended up.
[void]$ws.Range("I42").Select()
